$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new row 20 with the next forecast vector entry
# Clone the date cell's formatting (border/font/alignment/date numfmt) from the row above
$ws.Range("A19").Copy()
$ws.Range("A20").PasteSpecial(-4122)

$ws.Range("A20").Value = 45986

$ws.Range("B20").Value = 2025
$ws.Range("C20").Value = 2.622852459381209
$ws.Range("D20").Value = 2026
$ws.Range("E20").Value = 2.689750575689809
